$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Aufgabe 6 row (row 11): fill in the two empty cells ("Geschätzter" and
# "Tatsächlicher" Zeitlicher Aufwand) with "20". We rebuild each cell's
# paragraph via InsertXML so that the existing paragraph-mark run
# properties (e.g. the <w:lang w:val="de-DE"/> mark already present on the
# first cell) are carried over onto the newly created run, exactly as Word
# does when you type text into an empty, pre-formatted paragraph.

$cell2 = $t.Cell(11, 2)
$xmlCell2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="713D1D98" w14:textId="77777777" w:rsidR="002130F2" w:rsidRPr="00822462" w:rsidRDefault="002130F2" w:rsidP="00E14E24"><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>20</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cell2.Range.InsertXML($xmlCell2)

$t2 = $d.Tables.Item(1)
$cell3 = $t2.Cell(11, 3)
$xmlCell3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="50A544A9" w14:textId="77777777" w:rsidR="002130F2" w:rsidRDefault="002130F2" w:rsidP="00E14E24"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>20</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cell3.Range.InsertXML($xmlCell3)

# Aufgabe 10 (totals) row: Geschätzter Aufwand 60 -> 120, Tatsächlicher Aufwand 50 -> 105
$d.Content.Find.Execute("60", $true, $false, $false, $false, $false, $true, 1, $false, "120", 2)
$d.Content.Find.Execute("50", $true, $false, $false, $false, $false, $true, 1, $false, "105", 2)
